$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace lowercase purchase-invoice ids (i1..i10) in column A with
# uppercase (I1..I10)
$ws.Range("A2").Value = "I1"
$ws.Range("A3").Value = "I2"
$ws.Range("A4").Value = "I3"
$ws.Range("A5").Value = "I4"
$ws.Range("A6").Value = "I5"
$ws.Range("A7").Value = "I6"
$ws.Range("A8").Value = "I7"
$ws.Range("A9").Value = "I8"
$ws.Range("A10").Value = "I9"
$ws.Range("A11").Value = "I10"

# Replace lowercase shipment ids (s1..s4) in column B with uppercase
# (S1..S4)
$ws.Range("B2").Value = "S1"
$ws.Range("B3").Value = "S1"
$ws.Range("B4").Value = "S1"
$ws.Range("B5").Value = "S1"
$ws.Range("B6").Value = "S1"
$ws.Range("B7").Value = "S1"
$ws.Range("B8").Value = "S1"
$ws.Range("B9").Value = "S2"
$ws.Range("B10").Value = "S3"
$ws.Range("B11").Value = "S4"

# Update the view: drop the frozen/scrolled topLeftCell of C1 and move the
# active selection to B12
$ws.Range("B12").Select()
